$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (C column) date serial for all data rows (2-36) from 46070 to 46072
for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 3).Value = 46072
}

# Rows 13-30 got reordered (their Beteckning/Datum/Area values were shuffled between rows)
# Apply the new A (Beteckning), B (Datum), G (Area (ha)) values per row
$ws.Cells.Item(13, 1).Value = "A 21421-2021"
$ws.Cells.Item(13, 2).Value = 44316
$ws.Cells.Item(13, 7).Value = 0.6
$ws.Cells.Item(14, 1).Value = "A 28416-2024"
$ws.Cells.Item(14, 2).Value = 45477.62978009259
$ws.Cells.Item(14, 7).Value = 0.2
$ws.Cells.Item(15, 1).Value = "A 17491-2024"
$ws.Cells.Item(15, 2).Value = 45415.50266203703
$ws.Cells.Item(15, 7).Value = 6.2
$ws.Cells.Item(16, 1).Value = "A 50864-2022"
$ws.Cells.Item(16, 2).Value = 44867.56143518518
$ws.Cells.Item(16, 7).Value = 3.3
$ws.Cells.Item(17, 1).Value = "A 23503-2025"
$ws.Cells.Item(17, 2).Value = 45795
$ws.Cells.Item(17, 7).Value = 14.1
$ws.Cells.Item(18, 1).Value = "A 57001-2025"
$ws.Cells.Item(18, 2).Value = 45977
$ws.Cells.Item(18, 7).Value = 2
$ws.Cells.Item(19, 1).Value = "A 270-2025"
$ws.Cells.Item(19, 2).Value = 45660.48087962963
$ws.Cells.Item(19, 7).Value = 8.9
$ws.Cells.Item(20, 1).Value = "A 49633-2024"
$ws.Cells.Item(20, 2).Value = 45596.59559027778
$ws.Cells.Item(20, 7).Value = 0.8
$ws.Cells.Item(21, 1).Value = "A 24086-2025"
$ws.Cells.Item(21, 2).Value = 45795
$ws.Cells.Item(21, 7).Value = 0.7
$ws.Cells.Item(22, 1).Value = "A 4422-2024"
$ws.Cells.Item(22, 2).Value = 45327.45375
$ws.Cells.Item(22, 7).Value = 4.5
$ws.Cells.Item(23, 1).Value = "A 24212-2023"
$ws.Cells.Item(23, 2).Value = 45076
$ws.Cells.Item(23, 7).Value = 5.8
$ws.Cells.Item(24, 1).Value = "A 51434-2025"
$ws.Cells.Item(24, 2).Value = 45949
$ws.Cells.Item(24, 7).Value = 2.8
$ws.Cells.Item(25, 1).Value = "A 52965-2025"
$ws.Cells.Item(25, 2).Value = 45956
$ws.Cells.Item(25, 7).Value = 0.6
$ws.Cells.Item(26, 1).Value = "A 52960-2025"
$ws.Cells.Item(26, 2).Value = 45956
$ws.Cells.Item(26, 7).Value = 1.7
$ws.Cells.Item(27, 1).Value = "A 52888-2025"
$ws.Cells.Item(27, 2).Value = 45957.56943287037
$ws.Cells.Item(27, 7).Value = 2
$ws.Cells.Item(28, 1).Value = "A 50239-2022"
$ws.Cells.Item(28, 2).Value = 44865
$ws.Cells.Item(28, 7).Value = 13.2
$ws.Cells.Item(29, 1).Value = "A 46579-2024"
$ws.Cells.Item(29, 2).Value = 45582.75018518518
$ws.Cells.Item(29, 7).Value = 3
$ws.Cells.Item(30, 1).Value = "A 28409-2024"
$ws.Cells.Item(30, 2).Value = 45477.62280092593
$ws.Cells.Item(30, 7).Value = 0.4
